$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.079.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.35%  '
$ws.Range("D3").Value = "'1.623.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.06%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = "'214.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.98%  '
$ws.Range("E6").Value = '  -1.28%  '
$ws.Range("D8").Value = "'0.0633"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.01%  '
$ws.Range("E9").Value = '  -1.72%  '
$ws.Range("D10").Value = "'20.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.89%  '
$ws.Range("D11").Value = "'0.0848"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.09%  '
$ws.Range("D12").Value = "'1.850.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.04%  '
$ws.Range("D13").Value = "'1.646.17"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.19%  '
$ws.Range("E14").Value = '  +0.41%  '
$ws.Range("D15").Value = "'0.542"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'64.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.14%  '
$ws.Range("D17").Value = "'27.031.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.54%  '
$ws.Range("D18").Value = "'0.0₃0746"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.74%  '
$ws.Range("D19").Value = "'214.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.55%  '
$ws.Range("E20").Value = '  -0.13%  '
$ws.Range("E21").Value = '  -0.81%  '
$ws.Range("D22").Value = "'4.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.91%  '
$ws.Range("E23").Value = '  -6.62%  '
$ws.Range("E24").Value = '  -0.46%  '
$ws.Range("D25").Value = "'148.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.07%  '
$ws.Range("E26").Value = '  -0.20%  '
$ws.Range("D27").Value = "'7.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.58%  '
$ws.Range("E28").Value = '  -2.93%  '
$ws.Range("D29").Value = "'15.62"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.40%  '
$ws.Range("E30").Value = '  +0.78%  '
$ws.Range("E31").Value = '  -0.94%  '
$ws.Range("D32").Value = "'0.750"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +37.02%  '
$ws.Range("E33").Value = '  -0.84%  '
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("D35").Value = "'1.356.66"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.94%  '
$ws.Range("D36").Value = "'1.58"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.47%  '
$ws.Range("E37").Value = '  -0.74%  '
$ws.Range("E38").Value = '  +1.11%  '
$ws.Range("D39").Value = "'0.849"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.14%  '
$ws.Range("D41").Value = "'0.804"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.99%  '
$ws.Range("E42").Value = '  +0.15%  '
$ws.Range("D43").Value = "'65.16"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.35%  '
$ws.Range("E44").Value = '  +1.21%  '
$ws.Range("D45").Value = "'1.762.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.04%  '
$ws.Range("E46").Value = '  +32.27%  '
$ws.Range("D47").Value = "'90.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.62%  '
$ws.Range("E48").Value = '  +2.99%  '
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("D50").Value = "'0.102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.06%  '
$ws.Range("D51").Value = "'0.0514"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.44%  '
